# Alumni Fantasy League Results - apply "updated gitignore and added sos" edit
#
# Summary of changes:
#  1. NewData sheet ("Match Ups" table): append Week 3 .. Week 13 match-up
#     blocks (rows 14-79), each a merged/centered header row followed by five
#     match-up rows (Team, seed, seed, Team).
#  2. NewResults sheet (team stats table): column I ("Normalized WAE") switches
#     from text values ("0.8", "1.0", ...) to real numbers, and a new column J
#     "Strength of Schedule" is appended with win-count values.
#  3. View-state touch-ups: OldData tab is no longer the selected tab, the
#     Match Ups sheet's selection/scroll moves to the newly-added rows, and
#     NewResults becomes the active (displayed) sheet.

$wb = $excel.ActiveWorkbook

$oldData    = $wb.Worksheets.Item("OldData")
$matchUps   = $wb.Worksheets.Item("NewData")
$newResults = $wb.Worksheets.Item("NewResults")

# ---------------------------------------------------------------------------
# 1. Match Ups sheet: Week 3 - Week 13 blocks
# ---------------------------------------------------------------------------

$weeks = @(
    @{ Header = "Week 3"; Matchups = @(
            @("Tolosa",1,2,"Hallacy"),
            @("Nagle",3,4,"Stichler"),
            @("Rich",5,6,"Yamaoka"),
            @("Pitton",7,8,"Netter"),
            @("Walker",9,10,"Nishida")
    ) },
    @{ Header = "Week 4"; Matchups = @(
            @("Hallacy",1,2,"Rich"),
            @("Pitton",3,4,"Tolosa"),
            @("Nishida",5,6,"Yamaoka"),
            @("Nagle",7,8,"Netter"),
            @("Stichler",9,10,"Walker")
    ) },
    @{ Header = "Week 5"; Matchups = @(
            @("Nishida",1,2,"Hallacy"),
            @("Walker",3,4,"Nagle"),
            @("Yamaoka",5,6,"Stichler"),
            @("Netter",7,8,"Tolosa"),
            @("Rich",9,10,"Pitton")
    ) },
    @{ Header = "Week 6"; Matchups = @(
            @("Stichler",1,2,"Hallacy"),
            @("Tolosa",3,4,"Rich"),
            @("Nishida",5,6,"Pitton"),
            @("Walker",7,8,"Netter"),
            @("Nagle",9,10,"Yamaoka")
    ) },
    @{ Header = "Week 7"; Matchups = @(
            @("Hallacy",1,2,"Nagle"),
            @("Yamaoka",3,4,"Walker"),
            @("Stichler",5,6,"Pitton"),
            @("Netter",7,8,"Rich"),
            @("Nishida",9,10,"Tolosa")
    ) },
    @{ Header = "Week 8"; Matchups = @(
            @("Walker",1,2,"Hallacy"),
            @("Rich",3,4,"Nishida"),
            @("Nagle",5,6,"Pitton"),
            @("Yamaoka",7,8,"Netter"),
            @("Tolosa",9,10,"Stichler")
    ) },
    @{ Header = "Week 9"; Matchups = @(
            @("Hallacy",1,2,"Yamaoka"),
            @("Stichler",3,4,"Rich"),
            @("Nagle",5,6,"Tolosa"),
            @("Netter",7,8,"Nishida"),
            @("Pitton",9,10,"Walker")
    ) },
    @{ Header = "Week 10"; Matchups = @(
            @("Hallacy",1,2,"Netter"),
            @("Yamaoka",3,4,"Pitton"),
            @("Walker",5,6,"Tolosa"),
            @("Nishida",7,8,"Stichler"),
            @("Rich",9,10,"Nagle")
    ) },
    @{ Header = "Week 11"; Matchups = @(
            @("Hallacy",1,2,"Pitton"),
            @("Nagle",3,4,"Nishida"),
            @("Walker",5,6,"Rich"),
            @("Netter",7,8,"Stichler"),
            @("Tolosa",9,10,"Yamaoka")
    ) },
    @{ Header = "Week 12"; Matchups = @(
            @("Hallacy",1,2,"Tolosa"),
            @("Stichler",3,4,"Nagle"),
            @("Yamaoka",5,6,"Rich"),
            @("Pitton",7,8,"Netter"),
            @("Nishida",9,10,"Walker")
    ) },
    @{ Header = "Week 13"; Matchups = @(
            @("Rich",1,2,"Hallacy"),
            @("Tolosa",3,4,"Pitton"),
            @("Yamaoka",5,6,"Nishida"),
            @("Netter",7,8,"Nagle"),
            @("Walker",9,10,"Stichler")
    ) },
)

$row = 14
foreach ($week in $weeks) {
    # Header row: bold-free, centered "Week N" label merged across A:D
    $headerRange = $matchUps.Range("A" + $row + ":D" + $row)
    $matchUps.Cells.Item($row, 1).Value = $week.Header
    $matchUps.Cells.Item($row, 1).HorizontalAlignment = -4108   # xlCenter
    $headerRange.Merge()
    $row = $row + 1

    foreach ($m in $week.Matchups) {
        $matchUps.Cells.Item($row, 1).Value = $m[0]
        $matchUps.Cells.Item($row, 2).Value = $m[1]
        $matchUps.Cells.Item($row, 3).Value = $m[2]
        $matchUps.Cells.Item($row, 4).Value = $m[3]
        $row = $row + 1
    }
}

# ---------------------------------------------------------------------------
# 2. NewResults sheet: numeric "Normalized WAE" column + new "Strength of
#    Schedule" column
# ---------------------------------------------------------------------------

$newResults.Range("J1").Value = "Strength of Schedule"

$sosData = @(
    @(0.8, 4),
    @(1.0, 6),
    @(0.2, 9),
    @(1.1, 13),
    @(1.1, 9),
    @(1.3, 14),
    @(1.3, 8),
    @(1.1, 9),
    @(0.7, 13),
    @(1.3, 5)
)

$r = 2
foreach ($pair in $sosData) {
    $newResults.Cells.Item($r, 9).Value  = $pair[0]   # column I, now numeric
    $newResults.Cells.Item($r, 10).Value = $pair[1]   # column J
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3. View-state touch-ups
# ---------------------------------------------------------------------------

# Move selection/scroll on the Match Ups sheet down to the newly-added rows
$matchUps.Activate()
$matchUps.Range("B75:C79").Select()

# NewResults ends up as the active/displayed sheet
$newResults.Activate()
